# Actualización automática 2025-06-11 16:45:31
#
# Updates the monthly/group sales figures for "ARCOS GOMEZ CONSTRUCCIONES CIA. LTDA."
# (row 3 on both sheets) with a new sale of 57.86 registered for PORCELANATO /
# junio, and refreshes the dependent summary row (row 7) accordingly.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": new PORCELANATO sale for row 3 ---------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L3").Value = 57.86
$wsGrupo.Range("L7").Value = "1 de 5"

# --- Sheet "VENTA MENSUAL": new junio sale for row 3, and updated total ----
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F3").Value = 57.86
$wsMensual.Range("F7").Value = 57.86
